$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues

# Helper cell (outside the used range) used to stage text values so that
# Excel's automatic date recognition does not kick in when we write the
# new date-like strings ("12.08.11", "12.09.11", "12.10.11") into the sheet.
# We format the helper as Text, assign the literal string, copy it, and
# then paste only the *values* (not formats) into the real target cells,
# so the target cells keep their original (default) cell style.
$helper = $ws.Cells.Item(50, 50)

function Set-TextValue($cell, $text) {
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy()
    $cell.PasteSpecial($xlPasteValues)
}

# Rows 2-10: 12.06.11 -> 12.08.11 (col A), 12.06.11 -> 12.09.11 (col D)
for ($r = 2; $r -le 10; $r++) {
    Set-TextValue $ws.Cells.Item($r, 1) "12.08.11"
    Set-TextValue $ws.Cells.Item($r, 4) "12.09.11"
}

# Rows 11-18: 12.07.11 -> 12.09.11 (col A), 12.07.11 -> 12.10.11 (col D)
for ($r = 11; $r -le 18; $r++) {
    Set-TextValue $ws.Cells.Item($r, 1) "12.09.11"
    Set-TextValue $ws.Cells.Item($r, 4) "12.10.11"
}

# Clean up the helper cell so it leaves no trace in the used range
$helper.Clear()

# Update the selection to match the target state
$ws.Range("D12:D18").Select()
